# "Final training and evaluation" — final_model / model.pth are not
# committed for swin and deit, so this run's batch_size/negatives were
# halved (32->8, 16->4), early stopping was turned off, and the
# (start_run, end_run, best_val_loss) results for this sweep were
# cleared since the run never completed / was not committed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 6; $row++) {
    $ws.Range("J$row").Value  = 8      # batch_size: 32 -> 8
    $ws.Range("K$row").Value  = 4      # negatives:  16 -> 4
    $ws.Range("AE$row").Value = $false # use_early_stopping: TRUE -> FALSE

    # start_run / end_run / best_val_loss are no longer recorded
    $ws.Range("AF$row`:AH$row").ClearContents()
}

# Scroll the view over and select the (now empty) results columns, as in
# the saved workbook.
$ws.Range("AF2:AH6").Select()
